# Applies the 2025-12-05 CTA violent-crime YTD data refresh to before.xlsx
# Updates year-to-date cumulative counts across the Citywide Totals sheet,
# the By Neighborhood summary sheet, and the affected per-neighborhood sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 121
$ws.Range("K2").Value = 141
$ws.Range("E3").Value = 144
$ws.Range("J3").Value = 229
$ws.Range("L3").Value = 244
$ws.Range("B6").Value = 374
$ws.Range("C6").Value = 475
$ws.Range("D6").Value = 413
$ws.Range("E6").Value = 470
$ws.Range("F6").Value = 531
$ws.Range("I6").Value = 498
$ws.Range("J6").Value = 417
$ws.Range("B7").Value = 499
$ws.Range("C7").Value = 630
$ws.Range("D7").Value = 644
$ws.Range("E7").Value = 695
$ws.Range("F7").Value = 767
$ws.Range("I7").Value = 832
$ws.Range("J7").Value = 790
$ws.Range("K7").Value = 889
$ws.Range("L7").Value = 822

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("F7").Value = 11
$ws.Range("B8").Value = 31
$ws.Range("I21").Value = 15
$ws.Range("J29").Value = 13
$ws.Range("E32").Value = 66
$ws.Range("J32").Value = 47
$ws.Range("C36").Value = 39
$ws.Range("K50").Value = 30
$ws.Range("E53").Value = 82
$ws.Range("F53").Value = 82
$ws.Range("D65").Value = 25
$ws.Range("L74").Value = 10
$ws.Range("J78").Value = 6
$ws.Range("E91").Value = 7
$ws.Range("B96").Value = 16
$ws.Range("B98").Value = 499
$ws.Range("C98").Value = 630
$ws.Range("D98").Value = 644
$ws.Range("E98").Value = 695
$ws.Range("F98").Value = 767
$ws.Range("I98").Value = 832
$ws.Range("J98").Value = 790
$ws.Range("K98").Value = 889
$ws.Range("L98").Value = 822

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("F5").Value = 7
$ws.Range("F6").Value = 11

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("B6").Value = 22
$ws.Range("B7").Value = 31

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("I6").Value = 9
$ws.Range("I7").Value = 15

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J3").Value = 14
$ws.Range("E6").Value = 53
$ws.Range("J6").Value = 31
$ws.Range("E7").Value = 66
$ws.Range("J7").Value = 47

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("C6").Value = 34
$ws.Range("C7").Value = 39

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("K2").Value = 5
$ws.Range("K6").Value = 30

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("J4").Value = 5
$ws.Range("J5").Value = 6

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("B5").Value = 11
$ws.Range("B6").Value = 16

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("E3").Value = 13
$ws.Range("F6").Value = 61
$ws.Range("E7").Value = 82
$ws.Range("F7").Value = 82

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("E6").Value = 6
$ws.Range("E7").Value = 7

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("D5").Value = 24
$ws.Range("D6").Value = 25

$ws = $wb.Worksheets.Item('River North')
$ws.Range("L3").Value = 7
$ws.Range("L6").Value = 10

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("J2").Value = 1
$ws.Range("J6").Value = 13
